$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

# Corrected Cypher query text in B4: the "File Type" line was removed and the
# following line's indentation shifted by one extra space.
$query = @'
MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (f)-[*]->(c:case)<--(demo:demographic)
WHERE demo.breed IN ['Parson Russell Terrier'] 
OPTIONAL MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
WITH DISTINCT f, parent, c, demo, diag, s
RETURN  coalesce(f.file_name, '') AS `File Name`,
         coalesce(labels(parent)[0], '') AS `Association`,
        coalesce(f.file_description, '') AS `Description`,
        coalesce(f.file_format, '') AS `Format`,
        coalesce(f.file_size, '') AS `Size`,
        coalesce(c.case_id, '') AS `Case ID`,
        coalesce(demo.breed,'') AS Breed , 
        coalesce(diag.disease_term,'') AS Diagnosis , 
        coalesce(s.clinical_study_designation,'') AS `Study Code`
'@

$ws.Range("B4").Value = $query

# Row 4 shrinks now that it holds one fewer line of text.
$ws.Rows.Item(4).RowHeight = 232

# Selection moves from C12 to B4.
$ws.Range("B4").Select()
